# Update diaries/diary-ZihuaWeng.xlsx - add Jan 19th and Jan 20th diary entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Jan 19th entry (order chosen to match original authoring / shared-string order)
$ws.Range("A13").Value = "Jan 19th"
$ws.Range("B13").Value = "10pm-11.30pm"
$ws.Range("C13").Value = "Junxian, Wenchia"
$ws.Range("D13").Value = "Find a open source project for our coursework as our first choice was disapproved."
$ws.Range("E13").Value = "Checked out different Java open sourse project and learn to use IntellJ find the project size. At the end we chose Realm which is a mobile database that frequently used on Android "
$ws.Range("G13").Value = "excited!"
$ws.Range("F13").Value = "Many popular network framework used on Android do not have a lot code…Eleticsearch is perfect but too big for us..."
$ws.Rows.Item(13).RowHeight = 85

# Row 14: Jan 20th entry (order chosen to match original authoring / shared-string order)
$ws.Range("A14").Value = "Jan 20th"
$ws.Range("B14").Value = "5pm-8pm"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "Finish homework1"
$ws.Range("F14").Value = "Learn to commit and pull request after change multiple documents and after the fork repository got changed."
$ws.Range("G14").Value = "Happy to get my homework done!"
$ws.Range("E14").Value = "Finished homework1 and wrote a report using markdown. Refreshed my knowledge of how to write markdown document. Also get more understand of the project Jpacman3, especially the project structure and how the game wrote!"
$ws.Rows.Item(14).RowHeight = 119

# Update view/selection to reflect the edited area
$ws.Range("F13").Select()
